$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.772.30"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.69"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.88"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.16"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.70"
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.76"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.585.13"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.858"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.79"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.642.75"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.03"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.47"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.49"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.18"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.58"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.41"
$ws.Range("E29").Value = "  +9.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.26"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.72"
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.18"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0855"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.13"
$ws.Range("E35").Value = "  +5.13%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +8.50%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.120"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.84"
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.28"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.33"
$ws.Range("E41").Value = "  +24.92%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0320"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.802.10"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.201"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.28"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.69"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +7.32%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.54"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.18"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.28"
$ws.Range("E51").Value = "  +0.91%  "
